$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 5 (old rows 5-19 shift down to 6-20),
# copying formatting from the row above so date styling (s="2" on column D) carries over.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new weekly record.
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(5, 3).Value = "La Araucanía"
$ws.Cells.Item(5, 4).Value = 44811
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100108
$ws.Cells.Item(5, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(5, 9).Value = 100108001
$ws.Cells.Item(5, 10).Value = "Guayaba"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 60
$ws.Cells.Item(5, 14).Value = 2500
$ws.Cells.Item(5, 15).Value = 2500
$ws.Cells.Item(5, 16).Value = 2500
$ws.Cells.Item(5, 17).Value = "$/kilo"
$ws.Cells.Item(5, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 19).Value = 2500
$ws.Cells.Item(5, 20).Value = 1
